# Auto-generated Excel COM-interop script
# Applies the "test/train split (test=100)" change:
# 1. Deletes the "time" worksheet entirely (and its column/tab).
# 2. Updates the numeric coefficient values on every remaining sheet
#    to reflect the new test/train split recomputation.

$wb = $excel.ActiveWorkbook

# Step 1: remove the "time" sheet (deleted in the new split)
$timeSheet = $wb.Worksheets.Item("time")
$null = $timeSheet.Delete()

# Step 2: update the numeric values on every remaining sheet

$ws = $wb.Worksheets.Item("Gender")
$ws.Range("B2").Value = -0.0178163614334663
$ws.Range("C2").Value = -0.1354415532343072
$ws.Range("D2").Value = 0.04419580299056309
$ws.Range("E2").Value = 0.03766384537743038
$ws.Range("B3").Value = 0.01579717380434011
$ws.Range("C3").Value = 0.120091510534419
$ws.Range("D3").Value = -0.03918694531829928
$ws.Range("E3").Value = -0.03339527623465495

$ws = $wb.Worksheets.Item("Caste")
$ws.Range("B2").Value = -1.522489810969786
$ws.Range("C2").Value = 0.6644244153346088
$ws.Range("D2").Value = 0.2638110905242167
$ws.Range("E2").Value = 0.3930963071412782
$ws.Range("B3").Value = 0.5703723243055425
$ws.Range("C3").Value = -0.02895333672400313
$ws.Range("D3").Value = 0.317308244566617
$ws.Range("E3").Value = 0.231029386343876
$ws.Range("B4").Value = 1.727157275414552
$ws.Range("C4").Value = -0.559106336880007
$ws.Range("D4").Value = -0.1978052506118028
$ws.Range("E4").Value = -0.2987745830515177
$ws.Range("B5").Value = 3.1152185312969
$ws.Range("C5").Value = -1.814453376531952
$ws.Range("D5").Value = -1.242878005468966
$ws.Range("E5").Value = -1.475927810171124

$ws = $wb.Worksheets.Item("coaching")
$ws.Range("B2").Value = 0.37302176806933
$ws.Range("C2").Value = -0.2109162911484152
$ws.Range("D2").Value = 0.1137948932147942
$ws.Range("E2").Value = -0.1773209788383649
$ws.Range("B3").Value = -0.9388020756368896
$ws.Range("C3").Value = 0.8637560075195296
$ws.Range("D3").Value = 0.05812622744966152
$ws.Range("E3").Value = -0.3926038731240895
$ws.Range("B4").Value = 0.01237438203564128
$ws.Range("C4").Value = -0.0532136776215926
$ws.Range("D4").Value = -0.04405035513547918
$ws.Range("E4").Value = 0.1105683535882255

$ws = $wb.Worksheets.Item("Class_ten_education")
$ws.Range("B2").Value = -0.1207888431766439
$ws.Range("C2").Value = -0.05365959189744179
$ws.Range("D2").Value = -0.07855094008245525
$ws.Range("E2").Value = 0.1698199718799857
$ws.Range("B3").Value = -0.7486511925493202
$ws.Range("C3").Value = 0.4020675167857131
$ws.Range("D3").Value = -0.05575581210761441
$ws.Range("E3").Value = 0.03694317882397252
$ws.Range("B4").Value = 0.1172221339569855
$ws.Range("C4").Value = 0.00886036939752572
$ws.Range("D4").Value = 0.05087238205628844
$ws.Range("E4").Value = -0.1050640523051662

$ws = $wb.Worksheets.Item("twelve_education")
$ws.Range("B2").Value = -0.1678465912546628
$ws.Range("C2").Value = -0.05065792631240576
$ws.Range("D2").Value = -0.01617150903057662
$ws.Range("E2").Value = 0.1336665272389183
$ws.Range("B3").Value = 0.1921600301683246
$ws.Range("C3").Value = 0.04985608085337251
$ws.Range("D3").Value = 0.05518859675499277
$ws.Range("E3").Value = -0.1616710596207673
$ws.Range("B4").Value = 0.8170994419673605
$ws.Range("C4").Value = 0.4928416435255089
$ws.Range("D4").Value = -1.030680445130755
$ws.Range("E4").Value = -0.3892782724090617

$ws = $wb.Worksheets.Item("medium")
$ws.Range("B2").Value = 0.7417809481682555
$ws.Range("C2").Value = 0.2637815338315955
$ws.Range("D2").Value = -0.08783545354791973
$ws.Range("E2").Value = -0.7532915347975603
$ws.Range("B3").Value = 0.0455568163917259
$ws.Range("C3").Value = -0.05416500094834372
$ws.Range("D3").Value = -0.00391976547373472
$ws.Range("E3").Value = 0.05241628944331309
$ws.Range("B4").Value = -0.9759863302469386
$ws.Range("C4").Value = 0.1865369534347502
$ws.Range("D4").Value = 0.1043850236915534
$ws.Range("E4").Value = 0.2428076096328021

$ws = $wb.Worksheets.Item("Class_X_Percentage")
$ws.Range("B2").Value = -0.2710725833558689
$ws.Range("C2").Value = 0.4163957742273031
$ws.Range("D2").Value = 0.06256801988654964
$ws.Range("E2").Value = -0.224182188517642
$ws.Range("B3").Value = -0.08203039965249201
$ws.Range("C3").Value = 0.1095625979003029
$ws.Range("D3").Value = 0.01476941450493001
$ws.Range("E3").Value = 0.002500252151855765
$ws.Range("B4").Value = 0.007677152202905768
$ws.Range("C4").Value = 0.1994193185087823
$ws.Range("D4").Value = -0.05351654598344976
$ws.Range("E4").Value = -0.0673207604465511
$ws.Range("B5").Value = 0.4239550884900293
$ws.Range("C5").Value = -0.6445076300297451
$ws.Range("D5").Value = -0.0599040287002567
$ws.Range("E5").Value = 0.03601284192180274

$ws = $wb.Worksheets.Item("Class_XII_Percentage")
$ws.Range("B2").Value = -0.1173631909892116
$ws.Range("C2").Value = 0.4081734375564864
$ws.Range("D2").Value = 0.07356485366449239
$ws.Range("E2").Value = -0.456965389923419
$ws.Range("B3").Value = -0.160464320387934
$ws.Range("C3").Value = 0.1806052088994194
$ws.Range("D3").Value = -0.05455862713781379
$ws.Range("E3").Value = 0.0871748819317161
$ws.Range("B4").Value = 0.3898486193242168
$ws.Range("C4").Value = -0.1937842566689962
$ws.Range("D4").Value = 0.05003077063166095
$ws.Range("E4").Value = -0.2973418706972535
$ws.Range("B5").Value = 0.209465223388588
$ws.Range("C5").Value = -0.3556688166733264
$ws.Range("D5").Value = 0.09666691239817177
$ws.Range("E5").Value = -0.04121591754163066

$ws = $wb.Worksheets.Item("Father_occupation")
$ws.Range("B2").Value = -0.6127442099210322
$ws.Range("C2").Value = 0.4096266884755161
$ws.Range("D2").Value = -0.1544004407033456
$ws.Range("E2").Value = 0.1814592400730257
$ws.Range("B3").Value = -0.1382937410259995
$ws.Range("C3").Value = 0.07911685204837351
$ws.Range("D3").Value = 0.04424473682751962
$ws.Range("E3").Value = 0.01651397173452741
$ws.Range("B4").Value = -0.5601263269600919
$ws.Range("C4").Value = 0.5542456937151264
$ws.Range("D4").Value = -0.4728768016969819
$ws.Range("E4").Value = 0.2193520754319982
$ws.Range("B5").Value = -0.03927692677591278
$ws.Range("C5").Value = -0.08609944509297401
$ws.Range("D5").Value = 0.1984171680659208
$ws.Range("E5").Value = -0.1202511514211121
$ws.Range("B6").Value = -0.0519431066088205
$ws.Range("C6").Value = -0.08834832409243598
$ws.Range("D6").Value = 0.2043142601523672
$ws.Range("E6").Value = -0.1230423370853135
$ws.Range("B7").Value = 0.1045570764899687
$ws.Range("C7").Value = 0.09485025550905729
$ws.Range("D7").Value = -0.1191121680659469
$ws.Range("E7").Value = -0.004371230331547245
$ws.Range("B8").Value = 0.1285954612013405
$ws.Range("C8").Value = -0.04391166265906844
$ws.Range("D8").Value = -0.03264572778622188
$ws.Range("E8").Value = 0.006601007453734625
$ws.Range("B9").Value = 0.0752611468739896
$ws.Range("C9").Value = -0.1814357957314347
$ws.Range("D9").Value = 0.1151591541983272
$ws.Range("E9").Value = -0.04528434034438487

$ws = $wb.Worksheets.Item("Mother_occupation")
$ws.Range("B2").Value = -0.7276229805654599
$ws.Range("C2").Value = 0.246264038925953
$ws.Range("D2").Value = -1.624523595965311
$ws.Range("E2").Value = 0.9214280234331079
$ws.Range("B3").Value = -0.5048513266941949
$ws.Range("C3").Value = 0.1316452812892188
$ws.Range("D3").Value = -1.603845357632127
$ws.Range("E3").Value = 0.9415593110968651
$ws.Range("B4").Value = -0.4716525897304664
$ws.Range("C4").Value = 0.8861957008473945
$ws.Range("D4").Value = -1.11470670210302
$ws.Range("E4").Value = 0.2025438075549637
$ws.Range("B5").Value = -0.3768941077432657
$ws.Range("C5").Value = 0.7090702631245471
$ws.Range("D5").Value = -1.011966552438618
$ws.Range("E5").Value = 0.2970872541063498
$ws.Range("B6").Value = -0.2185312358448047
$ws.Range("C6").Value = 0.3642894971013477
$ws.Range("D6").Value = -0.5920458952175216
$ws.Range("E6").Value = 0.2865148388729123
$ws.Range("B7").Value = -0.6979472944825136
$ws.Range("C7").Value = -0.5960419866698066
$ws.Range("D7").Value = -0.934893794332917
$ws.Range("E7").Value = 1.208565539486198
$ws.Range("B8").Value = 0.02952048766629561
$ws.Range("C8").Value = -0.03847149773842611
$ws.Range("D8").Value = 0.0466204436774862
$ws.Range("E8").Value = 0.004018905409782753
$ws.Range("B9").Value = 0.07313893695974176
$ws.Range("C9").Value = 0.01399725327691484
$ws.Range("D9").Value = 0.211618763548545
$ws.Range("E9").Value = -0.2258889873095929
$ws.Range("B10").Value = 0.034490168156063
$ws.Range("C10").Value = -0.08900240772675676
$ws.Range("D10").Value = 0.1516913881702974
$ws.Range("E10").Value = -0.06712910679351529
